# Add a new "Merging Branches in a Local Repository" entry to column C,
# row 2, as rich text with three differently-formatted runs, matching the
# manual-steps style already used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$run1 = "Merging Branches in a Local RepositoryTo merge branches locally, use "
$run2 = "git checkout to switch to the branch you want to merge into"
$run3 = ". This branch is typically the main branch. Next, use git merge and specify the name of the other branch to bring into this branch."
$text = $run1 + $run2 + $run3

$c2 = $ws.Range("C2")
$c2.Value = $text

# Base/cell-level font for the whole string: Arial 10, color #1F1F1F.
$c2.Font.Name = "Arial"
$c2.Font.Size = 10
$c2.Font.Color = 2039583   # BGR(0x1F1F1F) -> RGB(0x1F,0x1F,0x1F)

# Second run ("git checkout ...") gets its own color, #040C28.
$start2 = $run1.Length + 1
$len2 = $run2.Length
$chars2 = $c2.Characters($start2, $len2)
$chars2.Font.Name = "Arial"
$chars2.Font.Size = 10
$chars2.Font.Color = 2624516   # RGB(0x04,0x0C,0x28)

# Third run (". This branch is typically ...") back to #1F1F1F.
$start3 = $start2 + $len2
$len3 = $run3.Length
$chars3 = $c2.Characters($start3, $len3)
$chars3.Font.Name = "Arial"
$chars3.Font.Size = 10
$chars3.Font.Color = 2039583   # RGB(0x1F,0x1F,0x1F)

# Move the active selection onto the new cell, as in the authored workbook.
$c2.Select() | Out-Null
